$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45901
$ws.Range("B2").Value = 45905

$ws.Range("L3").Select()
